# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps after regenerating the handback
# report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 794408b3-...md
$wsOverview.Range("G4").Value = "2016-08-21 02:49:54"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 794408b3-...xlf
$wsZhCn.Range("H4").Value = "2016-08-21 02:49:50"
$wsZhCn.Range("K4").Value = "2016-08-21 02:50:15"

# de-de sheet: Latest HO Xliff Generate Date (mirrors Overview) and
# Correspond Handback DateTime for 794408b3-...xlf
$wsDeDe.Range("H4").Value = "2016-08-21 02:49:54"
$wsDeDe.Range("K4").Value = "2016-08-21 02:50:22"
